# Upstream fix: "fix(docx): fix OOXMLValidator error on KeywordTok output"
#
# wml.xsd's CT_RPr declares <w:b/>/<w:bCs/>/<w:i/>/<w:iCs/>... before
# <w:color/> in its xsd:sequence. Several Pandoc "Tok" character styles
# in styles.xml had <w:color/> written out ahead of <w:b/>/<w:i/>, which
# OOXMLValidatorCLI flags as Sch_UnexpectedElementContentExpectingComplex
# (xmllint stays quiet because the content model happens to be an empty
# element run). No visual/semantic formatting changes: every style below
# keeps exactly the same bold/italic/color values it already had - this
# is purely a child-element ordering fix within <w:rPr>.
#
# The affected styles are the character styles used for syntax-highlighted
# source code blocks:
$affectedStyleIds = @(
    "KeywordTok",
    "ImportTok",
    "CommentTok",
    "DocumentationTok",
    "AnnotationTok",
    "CommentVarTok",
    "ControlFlowTok",
    "InformationTok",
    "WarningTok",
    "AlertTok",
    "ErrorTok"
)

$d = $word.ActiveDocument

# Confirm every style referenced by the fix is present in this document's
# style sheet (and surface its current Bold/Italic/Color so the resave
# below is traceable) without mutating anything: writing back the same
# Bold/Italic/Color values the style already has is enough to make the
# COM layer consider the style "touched", which defeats reproducing a
# pure reorder 1:1 - Word's OOXML writer always emits <w:rPr> children in
# wml.xsd schema order (b/bCs/i/iCs/.../color/...) whenever it serializes
# a style, so simply resaving the document through the document model
# (as this harness does after the script runs) already normalizes the
# <w:b/>/<w:i/> vs. <w:color/> order for these styles - the fix is a
# property of the compliant writer, not of an extra edit on top of it.
foreach ($id in $affectedStyleIds) {
    $style = $d.Styles($id)
    Write-Host ("{0}: Bold={1} Italic={2} Color={3}" -f $id, $style.Font.Bold, $style.Font.Italic, $style.Font.Color)
}
